$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: shift the Total row (20) and Footer row (21) down to rows 22-23
$ws.Rows.Item(20).Insert()
$ws.Rows.Item(20).Insert()

# Step 2: restore proper formatting on the two newly inserted blank rows (20-21)
#         by copying the formats of an existing, fully-styled item row (row 4)
$ws.Range("A4:N4").Copy()
$ws.Range("A20:N21").PasteSpecial(-4122)

# Step 3: merges for the two new item rows (matching the other item rows B:G / H:K / L:M)
$ws.Range("B20:G20").Merge()
$ws.Range("H20:K20").Merge()
$ws.Range("L20:M20").Merge()
$ws.Range("B21:G21").Merge()
$ws.Range("H21:K21").Merge()
$ws.Range("L21:M21").Merge()

# Step 4: write the full, alphabetically-sorted item list (18 rows) into rows 4-21,
#         including the 2 newly added products (ALKAPRESS, URIPAN)
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "ALKAPRESS PLUS 10/160MG 20 F.C. TABS."
$ws.Range("H4").Value = "0:0"
$ws.Range("L4").Value = 102
$ws.Range("N4").Value = "1:0"

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "DOLIPRANE 1 GM 15 TABS."
$ws.Range("H5").Value = "10:1"
$ws.Range("L5").Value = 48
$ws.Range("N5").Value = "1:0"

$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "ELICA-M CREAM 30 GRAM"
$ws.Range("H6").Value = "0:0"
$ws.Range("L6").Value = 52
$ws.Range("N6").Value = "1:0"

$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "GABAVERONA 300MG CAPS"
$ws.Range("H7").Value = "0:2"
$ws.Range("L7").Value = 42
$ws.Range("N7").Value = "0:0"

$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "KERELLA LOTION 30 ML"
$ws.Range("H8").Value = "3:0"
$ws.Range("L8").Value = 31
$ws.Range("N8").Value = "1:0"

$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "LIDOCAINE 10% TOPICAL SPRAY 15 GM"
$ws.Range("H9").Value = "1:0"
$ws.Range("L9").Value = 38
$ws.Range("N9").Value = "1:0"

$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "NEUROGLOPENTIN 300 MG 30 CAPS."
$ws.Range("H10").Value = "1:0"
$ws.Range("L10").Value = 37
$ws.Range("N10").Value = "0:0"

$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "SULBIN 750MG VIAL"
$ws.Range("H11").Value = "6:0"
$ws.Range("L11").Value = 35
$ws.Range("N11").Value = "1:0"

$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "SUPOLACK HAIR SHAMPOO 200 ML"
$ws.Range("H12").Value = "0:0"
$ws.Range("L12").Value = 149.5
$ws.Range("N12").Value = "1:0"

$ws.Range("A13").Value = 10
$ws.Range("B13").Value = "TELFAST 180MG 20 F.C. TABS"
$ws.Range("H13").Value = "1:0"
$ws.Range("L13").Value = 80
$ws.Range("N13").Value = "0:2"

$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "URIPAN 5MG 30 TAB."
$ws.Range("H14").Value = "1:1"
$ws.Range("L14").Value = 18
$ws.Range("N14").Value = "0:0"

$ws.Range("A15").Value = 12
$ws.Range("B15").Value = "URSOFALK 250MG 20 CAPS."
$ws.Range("H15").Value = "0:0"
$ws.Range("L15").Value = 122
$ws.Range("N15").Value = "1:0"

$ws.Range("A16").Value = 13
$ws.Range("B16").Value = "VIDROP 2800 I.U./ML ORAL DROPS 15 ML"
$ws.Range("H16").Value = "6:0"
$ws.Range("L16").Value = 26
$ws.Range("N16").Value = "1:0"

$ws.Range("A17").Value = 14
$ws.Range("B17").Value = "WELLMETAZONE 0.1% CREAM 40 GM"
$ws.Range("H17").Value = "0:0"
$ws.Range("L17").Value = 56
$ws.Range("N17").Value = "1:0"

$ws.Range("A18").Value = 15
$ws.Range("B18").Value = "WINZOXIB 90MG 20 TAB"
$ws.Range("H18").Value = "1:0"
$ws.Range("L18").Value = 66.64
$ws.Range("N18").Value = "0:0"

$ws.Range("A19").Value = 16
$ws.Range("B19").Value = "جنتيانا نقط"
$ws.Range("H19").Value = "4:0"
$ws.Range("L19").Value = 14
$ws.Range("N19").Value = "2:0"

$ws.Range("A20").Value = 17
$ws.Range("B20").Value = "سرنجات 5 سم"
$ws.Range("H20").Value = "-1:0"
$ws.Range("L20").Value = 2
$ws.Range("N20").Value = "1:0"

$ws.Range("A21").Value = 18
$ws.Range("B21").Value = "كريم فاتيكا 125 مل"
$ws.Range("H21").Value = "2:0"
$ws.Range("L21").Value = 50
$ws.Range("N21").Value = "1:0"

# Step 5: row heights for all data rows + total row + footer row
$ws.Rows.Item(4).RowHeight = 24.75
$ws.Rows.Item(5).RowHeight = 25.5
$ws.Rows.Item(6).RowHeight = 24.75
$ws.Rows.Item(7).RowHeight = 25.5
$ws.Rows.Item(8).RowHeight = 25.5
$ws.Rows.Item(9).RowHeight = 24.75
$ws.Rows.Item(10).RowHeight = 25.5
$ws.Rows.Item(11).RowHeight = 24.75
$ws.Rows.Item(12).RowHeight = 25.5
$ws.Rows.Item(13).RowHeight = 25.5
$ws.Rows.Item(14).RowHeight = 24.75
$ws.Rows.Item(15).RowHeight = 25.5
$ws.Rows.Item(16).RowHeight = 24.75
$ws.Rows.Item(17).RowHeight = 25.5
$ws.Rows.Item(18).RowHeight = 25.5
$ws.Rows.Item(19).RowHeight = 24.75
$ws.Rows.Item(20).RowHeight = 25.5
$ws.Rows.Item(21).RowHeight = 24.75
$ws.Rows.Item(22).RowHeight = 26.25
$ws.Rows.Item(23).RowHeight = 16.5

# Step 6: update the grand total (sum of price column) in the shifted Total row
$ws.Range("K22").Value = 969.14

Write-Output "edit complete"